$d = $word.ActiveDocument

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$wfull = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

# Process from the bottom of the document upward so paragraph indices
# earlier in the collection stay valid as we mutate later ones.

$paras = $d.Paragraphs

# 7) "Jumping on the general deals no damage" -> drop <w:lastRenderedPageBreak/>
$p62 = $paras.Item(62)
$xml62 = "<w:p $wfull w14:paraId='52076289' w14:textId='0DCB0454' w:rsidR='00302AEB' w:rsidRDefault='00302AEB' w:rsidP='00303772'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:cnfStyle w:val='000000000000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='0' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr><w:r><w:t>Jumping on the general deals no damage</w:t></w:r></w:p>"
$p62.Range.InsertXML($xml62)

# 6) "Can be defeated through powerups and normal attacks" -> drop <w:lastRenderedPageBreak/>
$p61 = $paras.Item(61)
$xml61 = "<w:p $wfull w14:paraId='1FDB22B4' w14:textId='77777777' w:rsidR='00F628D4' w:rsidRDefault='00302AEB' w:rsidP='00303772'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:cnfStyle w:val='000000000000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='0' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr><w:r><w:t>Can be defeated through powerups and normal attacks</w:t></w:r></w:p>"
$p61.Range.InsertXML($xml61)

# 5) "Punching the player deals a lot of damage" -> drop bookmarkStart/bookmarkEnd
$p60 = $paras.Item(60)
$xml60 = "<w:p $wfull w14:paraId='7DC6E335' w14:textId='69CC1424' w:rsidR='00B61D35' w:rsidRDefault='00B61D35' w:rsidP='00786273'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:cnfStyle w:val='000000000000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='0' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr><w:r><w:t>Punching the player deals a lot of damage</w:t></w:r></w:p>"
$p60.Range.InsertXML($xml60)

# 4) "Has more health than the mini General" -> drop <w:lastRenderedPageBreak/>
$p58 = $paras.Item(58)
$xml58 = "<w:p $wfull w14:paraId='74F98D96' w14:textId='7B22EE65' w:rsidR='00F628D4' w:rsidRDefault='00F628D4' w:rsidP='00786273'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:cnfStyle w:val='000000000000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='0' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr><w:r><w:t>Has more health than the mini General</w:t></w:r></w:p>"
$p58.Range.InsertXML($xml58)

# 3) "General" -> add <w:lastRenderedPageBreak/> before the text
$p55 = $paras.Item(55)
$xml55 = "<w:p $wfull w14:paraId='3E8B36EC' w14:textId='3479CAE7' w:rsidR='00F628D4' w:rsidRDefault='00F628D4' w:rsidP='00606AEB'><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>General</w:t></w:r></w:p>"
$p55.Range.InsertXML($xml55)

# 1 & 2) "Hat comes off when defeated" -> append extra run text, then insert
# a brand-new following paragraph with the extra bullet + bookmark.
$p53 = $paras.Item(53)
$xml53 = "<w:p $wfull w14:paraId='43B98E63' w14:textId='61F03B5A' w:rsidR='00772D67' w:rsidRDefault='003957A8' w:rsidP='00303772'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:cnfStyle w:val='000000100000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='1' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr><w:r><w:t>Hat comes off when defeated</w:t></w:r><w:r><w:t xml:space='preserve'> the first time</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:cnfStyle w:val='000000100000' w:firstRow='0' w:lastRow='0' w:firstColumn='0' w:lastColumn='0' w:oddVBand='0' w:evenVBand='0' w:oddHBand='1' w:evenHBand='0' w:firstRowFirstColumn='0' w:firstRowLastColumn='0' w:lastRowFirstColumn='0' w:lastRowLastColumn='0'/></w:pPr><w:r><w:t>When first &#8216;defeated&#8217; its health goes back up and player must defeat again</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
$p53.Range.InsertXML($xml53)
